$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update TAL ferm improvements scenario values
$ws.Range("E23").Value = 0.73
$ws.Range("E24").Value = 68

# Minor plotting change: update the active selection
$ws.Range("E20").Select()
